# removed ALL option in supplier scope for #% of Vol, fixed logic for # of volume awarded
#
# This script:
#  1. Sorts the "Price" sheet data (A1:C31) ascending by column B (Bid ID / Facility),
#     then removes the "ALL"-scope row for Supplier C / Bid 1, zeroes out the price
#     for Supplier A / Bid 1, and blanks the price for Supplier B / Bid 1.
#  2. Sorts the "Supplier Bid Attributes" sheet data (A1:E31) ascending by column B.
#  3. Re-applies AutoFilter to both sheets (and registers the corresponding
#     hidden _xlnm._FilterDatabase workbook-level names).
#  4. Switches the active sheet/selection from "Supplier Bid Attributes" to "Price".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Price sheet
# ---------------------------------------------------------------------------
$wsPrice = $wb.Worksheets.Item("Price")

$priceRange = $wsPrice.Range("A1:C31")
$priceKey1 = $wsPrice.Range("B2")
$priceRange.Sort($priceKey1, 1, $null, $null, 1, $null, $null, 1)

# After the sort, row 4 is Supplier "C" / Bid "1" (the removed "ALL" scope entry).
$wsPrice.Rows.Item(4).Delete()

# Supplier "A" / Bid "1" -> Price 0 ; Supplier "B" / Bid "1" -> Price cleared.
$wsPrice.Cells.Item(2, 3).Value = 0
$wsPrice.Cells.Item(3, 3).ClearContents()

$priceFilterRange = $wsPrice.Range("A1:C1")
$priceFilterRange.AutoFilter()
$priceName = $wsPrice.Names.Add("_xlnm._FilterDatabase", "=Price!`$A`$1:`$C`$1")
$priceName.Visible = $false

# ---------------------------------------------------------------------------
# 2. Supplier Bid Attributes sheet
# ---------------------------------------------------------------------------
$wsSupplier = $wb.Worksheets.Item("Supplier Bid Attributes")

$supplierRange = $wsSupplier.Range("A1:E31")
$supplierKey1 = $wsSupplier.Range("B2")
$supplierRange.Sort($supplierKey1, 1, $null, $null, 1, $null, $null, 1)

$supplierFilterRange = $wsSupplier.Range("A1:E1")
$supplierFilterRange.AutoFilter()
$supplierName = $wsSupplier.Names.Add("_xlnm._FilterDatabase", "=`'Supplier Bid Attributes`'!`$A`$1:`$E`$1")
$supplierName.Visible = $false

# ---------------------------------------------------------------------------
# 3. Active sheet / selection bookkeeping
#    (select the previously-active sheet first so the final Activate/Select
#    on "Price" is what ends up recorded as the active tab & selection)
# ---------------------------------------------------------------------------
$wsSupplier.Activate()
$wsSupplier.Range("C7").Select()

$wsPrice.Activate()
$wsPrice.Range("I10").Select()
